$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C87").Value = 1608
$ws.Range("C88").Value = 1635
$ws.Range("C89").Value = 1664
$ws.Range("C90").Value = 1684
$ws.Range("C91").Value = 1708
$ws.Range("C92").Value = 1730

$ws.Range("D90").Formula = "=C90-C89"
$ws.Range("D91").Formula = "=C91-C90"
$ws.Range("D92").Formula = "=C92-C91"

$ws.Range("E90").Formula = "=(C90-C83)/7"
$ws.Range("E91").Formula = "=(C91-C84)/7"
$ws.Range("E92").Formula = "=(C92-C85)/7"

$ws.Range("A65").Select()
$ws.Range("V78").Select()
